$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData_getToken")

# --- New "Status Code" column (D) data, filled first so the new shared
#     strings ("200"/"500") land before the renamed header strings, matching
#     the order the workbook was actually authored in. ---
$ws.Range("D2").Value = "200"
$ws.Range("D3").Value = "500"
$ws.Range("D4").Value = "500"
$ws.Range("D5").Value = "500"
$ws.Range("D6").Value = "200"

# --- Header row: rename + add the Status Code header ---
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Status Code"

# --- Column D width (displays as width 11 with best-fit in the saved file) ---
$ws.Columns.Item(4).ColumnWidth = 10.2

# --- Drop the border that used to ring the TC_No column's data cells ---
$ws.Range("A2:A6").Borders.LineStyle = -4142

# --- Match the saved selection/active cell ---
$ws.Range("E2").Select() | Out-Null
